# Apply "contingencies with rene fine" edit:
# - Insert two new line entries (line7, line8) taking over rows 8 & 9
#   (pushing the existing extr1..extr8 rows down by two positions).
# - Re-populate C (from_bus), D (to_bus) and E (in_service) columns with
#   their new values for every affected row.
# - Append two brand-new rows (16 & 17) for extr7 / extr8 so all eight
#   extr* entries are still present after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create rows 16 & 17 by cloning the formatting of row 15 (A column
#     style, borders, bold, etc.) before filling in values. ---
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A15:E15").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row-by-row target data (row, A, B-name, C, D, E) ---
$rows = @(
    @{ Row = 8;  A = 6;  Name = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  Name = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  Name = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  Name = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; Name = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; Name = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; Name = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; A = 13; Name = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; Name = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; Name = "extr8"; C = 8;  D = 5;  E = $true  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.Name
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
}
